$d = $word.ActiveDocument

$pairs = @(
    @("765÷5=", "274÷6="),
    @("605÷7=", "187÷6="),
    @("893÷8=", "418÷7="),
    @("530÷8=", "700÷8="),
    @("602÷8=", "780÷7="),
    @("110÷9=", "343÷3="),
    @("888÷2=", "336÷7="),
    @("285÷9=", "418÷7="),
    @("871÷2=", "655÷4="),
    @("513÷7=", "290÷5="),
    @("832÷4=", "682÷6="),
    @("793÷3=", "544÷7="),
    @("303÷4=", "590÷3="),
    @("990÷8=", "679÷4="),
    @("117÷2=", "262÷9="),
    @("769÷3=", "821÷3="),
    @("643÷9=", "119÷5="),
    @("991÷6=", "881÷7="),
    @("613÷7=", "652÷8="),
    @("100÷7=", "957÷3="),
    @("191÷5=", "678÷4="),
    @("370÷4=", "395÷6="),
    @("869÷4=", "515÷8="),
    @("995÷8=", "535÷5="),
    @("567÷8=", "625÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
